$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45436
$ws.Range("D33").Value = 576.08
$ws.Range("D34").Value = 630.0890000000001
$ws.Range("D35").Value = 706.599
$ws.Range("D36").Value = 859.62

# Re-merge the cells so they are re-registered in the sheet's merge list
# in the same order Excel wrote them after the edit.
$mergeRanges = @("B33:C33", "B34:C34", "B36:C36", "B32:C32", "A11:D11", "A10:D10", "B35:C35")
foreach ($ref in $mergeRanges) {
    $rng = $ws.Range($ref)
    $rng.MergeCells = $false
    $rng.Merge()
}
